$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates (sound volume / difficulty / enemies / controls additions) ---

# Week 2 plan row (row 5, "détail" row for Etudiant 2)
$ws.Range("C5").Value = "40min prise en main git`nmaquette + répartition fonctionnalités`nRédactions des règles /but`n2h Classe Ennemis"
$ws.Range("D5").Value = "2h fonctionnalitées`nde la classe ennemi"
$ws.Range("D5").WrapText = $true
$ws.Range("D5").VerticalAlignment = -4160
$ws.Range("E5").Value = "Pathfinding`nCode de triche (Joueur & Debug)`nCôntroles d'utilisateur"
$ws.Range("G5").Value = "Generation Ennemis s/map`nCollisions ennemis et`njoueur"

# Week 2 second block (row 12, "détail" row for Etudiant 2)
$ws.Range("C12").Value = "`nSelection Difficulté`nBarre de vie`n"
$ws.Range("D12").Value = "Menu son, parametres`nfix  etc.."
$ws.Range("D12").WrapText = $true
$ws.Range("D12").VerticalAlignment = -4160

# --- View state: selection moved to B10, scrolled down to show row 8 ---
$ws.Range("B10").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
